$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear formulas/values in A32:B36, leaving the cells blank (style retained)
$ws.Range("A32:B36").ClearContents()

# Update the view: scroll so row 22 is at top, and select A32:B36 with active cell B36
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("A32:B36").Select()
$ws.Range("B36").Activate()
